$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col_B = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
$col_C = @(1.028272275890204, 1.029492143617646, 1.030281606163354, 1.030613528437029, 1.030669261592588, 1.030286041196023, 1.028684510851954, 1.025863292912474, 1.023982952377646, 1.023168828136641, 1.022866435022286, 1.022931299009981, 1.023143832060746, 1.023274781789852, 1.024036984503565, 1.024515112621055, 1.024794004127799, 1.024889100265232, 1.024463813246012, 1.023081246149432, 1.02221202147391, 1.02267280999887, 1.024486993202379, 1.026592553311649)
$col_D = @(1.033364396552964, 1.033770693757986, 1.034033190620537, 1.034143445376686, 1.034161951774116, 1.034034664241319, 1.033501789133646, 1.032559787080287, 1.031929885495124, 1.031656706225684, 1.031555173105067, 1.031576955084414, 1.031648314722646, 1.031692273574002, 1.031948006670106, 1.032108308246415, 1.032201768035729, 1.032233628321391, 1.032091113672518, 1.031627302787493, 1.031335327761509, 1.031490142557016, 1.032098883288621, 1.032803661011443)
$col_E = @(1.037472196307419, 1.038556118407225, 1.039257794178608, 1.039552852725858, 1.039602398722026, 1.03926173647331, 1.037838450084049, 1.035332739492436, 1.033663755719051, 1.032941407047278, 1.032673143240644, 1.032730684556675, 1.032919231292222, 1.033035407567584, 1.033711702514256, 1.034136012538422, 1.034383537547902, 1.034467942670902, 1.034090484773108, 1.032863707645208, 1.032092665029592, 1.03250138311362, 1.034111056711302, 1.035980258885621)
$col_F = @(1.046638173760821, 1.047809647248494, 1.048567726082159, 1.048886436926925, 1.048939950754449, 1.048571984654083, 1.04703406737161, 1.044324461043154, 1.042518270742739, 1.041736207849456, 1.041445718029648, 1.04150802895938, 1.041712195813178, 1.041837990131103, 1.042570174193771, 1.043029460901174, 1.043297357993076, 1.043388704512142, 1.04298018349167, 1.041652073686374, 1.04081705599475, 1.041259713421169, 1.043002449822994, 1.045024917453846)
$col_I = @(1.032622301507724, 1.032677442911459, 1.032711461240549, 1.032725363918995, 1.032727674841483, 1.032711648576138, 1.032641280433504, 1.032504589134172, 1.032404972766827, 1.032359833850643, 1.032342766851698, 1.032346441362643, 1.032358429210915, 1.032365775536542, 1.032407926321386, 1.032433830290818, 1.032448746196096, 1.032453799300614, 1.032431071043435, 1.032354907374874, 1.032305282618539, 1.03233175410334, 1.032432318426394, 1.03254142664946)
$col_J = @(1.033425514789117, 1.034284535295681, 1.034840028273569, 1.035073474524338, 1.035112666338216, 1.034843147918725, 1.033715897995621, 1.031726809762031, 1.030398852442587, 1.029823367271793, 1.029609534360386, 1.029655405529934, 1.029805693255176, 1.029898280889036, 1.030437035461833, 1.030774854948941, 1.030971853998731, 1.031039017933477, 1.03073861481558, 1.029761439266959, 1.029146632422734, 1.029472593217048, 1.030754990312062, 1.032241366737783)
$col_K = @(1.036166958111639, 1.036383295965178, 1.036522139374415, 1.036580234418838, 1.036589972681188, 1.036522916724518, 1.03624030604183, 1.03573362826684, 1.035390089417513, 1.035239986128721, 1.035184029856512, 1.035196041750449, 1.035235364862727, 1.035259566507126, 1.035400022982576, 1.035487767524604, 1.035538817208417, 1.035556201700441, 1.035478366825275, 1.035223790729937, 1.035062565080279, 1.035148143674574, 1.035482614999258, 1.035865636807331)
$col_L = @(1.040262983439837, 1.041155994684157, 1.041733548070532, 1.041976284186464, 1.042017036723727, 1.041736791785957, 1.040564840185884, 1.038497495497203, 1.03711772223817, 1.036519886474246, 1.036297764656489, 1.036345413229086, 1.03650152702584, 1.036597706035504, 1.037157390420267, 1.037508361784532, 1.037713040392968, 1.037782824259876, 1.037470709720396, 1.036455557041804, 1.035816948340192, 1.036155519700782, 1.037487723182027, 1.039032221774529)
$col_M = @(1.049403049437539, 1.050385271528423, 1.051020318311885, 1.051287168954063, 1.051331967155563, 1.051023884465858, 1.049735103345248, 1.047460120360787, 1.045940740825306, 1.045282175781077, 1.045037454647293, 1.045089952737337, 1.045261949117217, 1.045367908381309, 1.045984433546451, 1.046370985072707, 1.046596390074652, 1.046673236537018, 1.046329518376774, 1.045211303259939, 1.044507653583812, 1.04488072724492, 1.046348255589508, 1.048048734893587)
$col_N = @(1.034893097555681, 1.035753337969895, 1.036309619811528, 1.036543397582749, 1.036582645053499, 1.036312743886937, 1.035183893139645, 1.033191980172612, 1.031862137001455, 1.031285834575711, 1.031071697997042, 1.031117634308912, 1.031268135459962, 1.031360854578882, 1.031900374244969, 1.032238673474497, 1.032435952285534, 1.032503211600771, 1.032202381875989, 1.031223818626006, 1.030608138685543, 1.030934562381573, 1.032218780627555, 1.03370726787829)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $col_B[$i]
    $ws.Range("C$row").Value = $col_C[$i]
    $ws.Range("D$row").Value = $col_D[$i]
    $ws.Range("E$row").Value = $col_E[$i]
    $ws.Range("F$row").Value = $col_F[$i]
    $ws.Range("I$row").Value = $col_I[$i]
    $ws.Range("J$row").Value = $col_J[$i]
    $ws.Range("K$row").Value = $col_K[$i]
    $ws.Range("L$row").Value = $col_L[$i]
    $ws.Range("M$row").Value = $col_M[$i]
    $ws.Range("N$row").Value = $col_N[$i]
}
